$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# F32: "not entered" -> "entered"  (only existing cell whose text really
# changes; every other existing s="..." cell only shifts shared-string index
# because new strings are inserted earlier in the table, which the engine
# re-resolves automatically once we stop referencing the old text anywhere)
# ---------------------------------------------------------------------------
$ws.Range("F32").Value2 = "entered"

# ---------------------------------------------------------------------------
# New data rows 33-38 (CAISN DFO campaign GPS points), written in natural
# row-major order so shared strings get appended in the same order Excel
# itself would assign them on save.
# ---------------------------------------------------------------------------

# Row 33 - Deception Bay 4
$ws.Range("A33").Value2 = "CAISN DFO"
$ws.Range("B33").Value2 = 41122
$ws.Range("C33").Value2 = "Deception Bay 4"
$ws.Range("D33").Value2 = 62.24567
$ws.Range("E33").Value2 = -74.86087
$ws.Range("F33").Value2 = "not entered"

# Row 34 - Steensby Inlet T1
$ws.Range("A34").Value2 = "CAISN DFO"
$ws.Range("B34").Value2 = 41136
$ws.Range("C34").Value2 = "Steensby Inlet T1"
$ws.Range("D34").Value2 = 62.24567
$ws.Range("E34").Value2 = -74.86087
$ws.Range("F34").Value2 = "not entered"

# Row 35 - Steensby Inlet T2
$ws.Range("A35").Value2 = "CAISN DFO"
$ws.Range("B35").Value2 = 41140
$ws.Range("C35").Value2 = "Steensby Inlet T2"
$ws.Range("D35").Value2 = 70.21439
$ws.Range("E35").Value2 = -78.76381
$ws.Range("F35").Value2 = "not entered"

# Row 36 - Steensby Inlet T3
$ws.Range("A36").Value2 = "CAISN DFO"
$ws.Range("B36").Value2 = 41138
$ws.Range("C36").Value2 = "Steensby Inlet T3"
$ws.Range("D36").Value2 = 70.19653
$ws.Range("E36").Value2 = -78.39088
$ws.Range("F36").Value2 = "not entered"

# Row 37 - Steensby Inlet T5
$ws.Range("A37").Value2 = "CAISN DFO"
$ws.Range("B37").Value2 = 41139
$ws.Range("C37").Value2 = "Steensby Inlet T5"
$ws.Range("D37").Value2 = 70.27333
$ws.Range("E37").Value2 = -78.51007
$ws.Range("F37").Value2 = "not entered"

# Row 38 - Steensby Inlet T6
$ws.Range("A38").Value2 = "CAISN DFO"
$ws.Range("B38").Value2 = 41139
$ws.Range("C38").Value2 = "Steensby Inlet T6"
$ws.Range("D38").Value2 = 70.32535
$ws.Range("E38").Value2 = -78.56289
$ws.Range("F38").Value2 = "not entered"

# ---------------------------------------------------------------------------
# Formatting to approximate the hand-pasted look of the new rows:
#  - date columns (B33:B38) -> same date format as the rest of column B
#  - lat/long columns get a dedicated font + centred alignment (two slightly
#    different fonts were used when this block was pasted in, Tahoma for the
#    Steensby Inlet rows and Calibri for the Deception Bay / T1 rows)
# ---------------------------------------------------------------------------

# Reuse the workbook's existing date format (numFmtId 14) for column B.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B33:B38").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Deception Bay 4 / Steensby Inlet T1 lat-long values: Calibri, centred.
$llRange1 = $ws.Range("D33:E34")
$llRange1.Font.Name = "Calibri"
$llRange1.Font.Size = 11
$llRange1.HorizontalAlignment = -4108

# Steensby Inlet T2/T3/T5/T6 lat-long values: Tahoma, centred.
$llRange2 = $ws.Range("D35:E38")
$llRange2.Font.Name = "Tahoma"
$llRange2.Font.Size = 10
$llRange2.HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# Sheet view: scroll back to the top and move the selection to G7 (matches
# the saved cursor position recorded in the workbook).
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G7").Select() | Out-Null

# ---------------------------------------------------------------------------
# Page setup: the sheet was switched to explicit portrait orientation.
# ---------------------------------------------------------------------------
$ws.PageSetup.Orientation = 1
